$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.287.58'
$ws.Range('E2').Value = '  +0.05%  '

$ws.Range('D3').Value = '1.917.84'
$ws.Range('E3').Value = '  -0.58%  '

$ws.Range('D4').Value = "'" + '1.013'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +1.42%  '

$ws.Range('D5').Value = "'" + '0.7327'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.77%  '

$ws.Range('D6').Value = "'" + '242.39'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.85%  '

$ws.Range('D7').Value = "'" + '1.006'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.80%  '

$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Value = "'" + '0.3105'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.95%  '

$ws.Range('B9').Value = 'Solana'
$ws.Range('C9').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D9').Value = "'" + '27.01'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.38%  '

$ws.Range('D10').Value = "'" + '0.06908'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.74%  '

$ws.Range('D11').Value = "'" + '0.08021'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.02%  '

$ws.Range('D12').Value = "'" + '0.7628'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.37%  '

$ws.Range('D13').Value = '1.921.40'
$ws.Range('E13').Value = '  -0.36%  '

$ws.Range('D14').Value = "'" + '5.290'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.78%  '

$ws.Range('D15').Value = "'" + '91.11'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.78%  '

$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '30.204.74'
$ws.Range('E16').Value = '  -0.17%  '

$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').Value = "'" + '14.09'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -3.70%  '

$ws.Range('D18').Value = "'" + '245.76'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -3.58%  '

$ws.Range('D19').Value = "'" + '5.813'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.13%  '

$ws.Range('D20').Value = "'" + '0.000007758'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.74%  '

$ws.Range('D21').Value = "'" + '1.007'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.88%  '

$ws.Range('D22').Value = '2.147.92'
$ws.Range('E22').Value = '  -1.56%  '

$ws.Range('D23').Value = "'" + '1.011'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.97%  '

$ws.Range('D24').Value = "'" + '6.562'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.89%  '

$ws.Range('D25').Value = "'" + '9.375'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.11%  '

$ws.Range('D26').Value = "'" + '165.03'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.54%  '

$ws.Range('D27').Value = "'" + '18.86'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.27%  '

$ws.Range('D28').Value = "'" + '0.1270'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -5.78%  '

$ws.Range('D29').Value = "'" + '2.130'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -7.37%  '

$ws.Range('D30').Value = "'" + '1.374'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.38%  '

$ws.Range('D31').Value = "'" + '1.549'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.18%  '

$ws.Range('D32').Value = "'" + '4.333'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.41%  '

$ws.Range('D33').Value = "'" + '4.036'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.76%  '

$ws.Range('D34').Value = "'" + '0.05158'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.82%  '

$ws.Range('D35').Value = "'" + '1.283'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.63%  '

$ws.Range('D36').Value = "'" + '0.7436'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.75%  '

$ws.Range('D37').Value = "'" + '2.788'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.76%  '

$ws.Range('D38').Value = "'" + '0.01918'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.96%  '

$ws.Range('D39').Value = "'" + '2.778'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.71%  '

$ws.Range('D40').Value = "'" + '6.442'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.58%  '

$ws.Range('D41').Value = "'" + '75.60'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.31%  '

$ws.Range('D42').Value = "'" + '0.4416'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.17%  '

$ws.Range('D43').Value = "'" + '1.919'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.41%  '

$ws.Range('D44').Value = "'" + '1.006'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.75%  '

$ws.Range('D45').Value = "'" + '0.8377'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.47%  '

$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = "'" + '100.74'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.80%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'" + '9.812'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.11%  '

$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = "'" + '7.533'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.08%  '

$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.055.37'
$ws.Range('E49').Value = '  -1.86%  '

$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = "'" + '36.88'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.24%  '

$ws.Range('D51').Value = "'" + '0.1199'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.70%  '
